# Added Diode, Began Routing
#
# Adds a new part row (row 19) to the "Comprehensive Parts List" sheet
# describing a regulator/USB protection diode, and extends the existing
# shared "price comp" formula down to cover it.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Comprehensive Parts List")

$ws.Range("A19").Value = "Regulator Diode"
$ws.Range("B19").Value = "CTS05S40L3FCT-ND"
$ws.Range("C19").Value = "SOD-882"
$ws.Range("D19").Value = 0.35
$ws.Range("E19").Value = 1
$ws.Range("F19").Value = "Diode between USB and VIN-(Unregulated)"
$ws.Range("G19").Formula = "=D19*E19"

$ws.Activate()
$ws.Range("C20").Select()
